# The "exponential" sheet (first sheet) lists model results in columns B:E,
# each column holding a header label in row 1 and a value in row 2.
# This change reverses the order of those four columns
# (PFS caba, PFS mito, OS caba, OS mito) -> (OS mito, OS caba, PFS mito, PFS caba)
# while keeping each label together with its corresponding value.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("exponential")

$headers = @()
$values = @()
foreach ($col in @("B", "C", "D", "E")) {
    $headers += $ws.Range("$col" + "1").Value()
    $values += $ws.Range("$col" + "2").Value()
}

$newCols = @("B", "C", "D", "E")
$srcOrder = @(3, 2, 1, 0)  # reversed order

for ($i = 0; $i -lt 4; $i++) {
    $col = $newCols[$i]
    $src = $srcOrder[$i]
    $ws.Range("$col" + "1").Value = $headers[$src]
    $ws.Range("$col" + "2").Value = $values[$src]
}
